$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.293.31"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.548.54"
$ws.Range("E3").Value = "  -2.31%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.68"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.63"
$ws.Range("E6").Value = "  +5.62%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -0.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.547.85"
$ws.Range("E9").Value = "  -2.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.139"
$ws.Range("E10").Value = "  +1.51%  "
$ws.Range("E11").Value = "  +1.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.352"
$ws.Range("E12").Value = "  -2.69%  "
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.03"
$ws.Range("E14").Value = "  -0.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.013.12"
$ws.Range("E15").Value = "  -2.27%  "
$ws.Range("E16").Value = "  -0.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.161.46"
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.544.57"
$ws.Range("E18").Value = "  -2.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.06"
$ws.Range("E19").Value = "  +3.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.40"
$ws.Range("E20").Value = "  -3.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "355.27"
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.22"
$ws.Range("E22").Value = "  -1.07%  "
$ws.Range("E23").Value = "  +1.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.99"
$ws.Range("E24").Value = "  +3.80%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.11"
$ws.Range("E26").Value = "  +1.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.01"
$ws.Range("E27").Value = "  -4.96%  "
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0998"
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "535.87"
$ws.Range("E31").Value = "  -0.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.25"
$ws.Range("E32").Value = "  +5.02%  "
$ws.Range("E33").Value = "  +0.84%  "
$ws.Range("E34").Value = "  -0.54%  "
$ws.Range("E35").Value = "  -1.14%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  -0.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "157.60"
$ws.Range("E38").Value = "  +0.32%  "
$ws.Range("E39").Value = "  -0.45%  "
$ws.Range("E40").Value = "  +1.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.357"
$ws.Range("E41").Value = "  -1.76%  "
$ws.Range("E42").Value = "  +0.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.19"
$ws.Range("E43").Value = "  +1.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.57"
$ws.Range("E44").Value = "  +7.10%  "
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.83"
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "151.89"
$ws.Range("E47").Value = "  +0.56%  "
$ws.Range("E48").Value = "  -1.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0280"
$ws.Range("E49").Value = "  -6.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.72"
$ws.Range("E50").Value = "  -1.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.72"
$ws.Range("E51").Value = "  +1.63%  "
